$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("course")

# New course rows (course_offered endpoints added for CMPE/STAT courses)
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = 2
$ws.Cells.Item(33,3).Value = "Principles of Digital Design"
$ws.Cells.Item(33,4).Value = 212
$ws.Cells.Item(33,5).Value = 4

$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = 2
$ws.Cells.Item(34,3).Value = "Systems Design and Programming"
$ws.Cells.Item(34,4).Value = 310
$ws.Cells.Item(34,5).Value = 4

$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = 2
$ws.Cells.Item(35,3).Value = "Probability, Statistics, and Random Processes"
$ws.Cells.Item(35,4).Value = 320
$ws.Cells.Item(35,5).Value = 3

$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = 4
$ws.Cells.Item(36,3).Value = "Introduction to Probability and Statistics for Scientists and Engineers"
$ws.Cells.Item(36,4).Value = 355
$ws.Cells.Item(36,5).Value = 4

# Make "course" the active sheet/tab and move its selection to the newly added data
$ws.Activate()
$ws.Range("E37").Select()
